$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 112
$ws.Range("I5").Value = 129.35715
$ws.Range("J5").Value = 63.4
$ws.Range("K5").Value = 129.35715
$ws.Range("L5").Value = 63.4
$ws.Range("M5").Value = -14.35714999999999
$ws.Range("N5").Value = -293.4

$ws.Range("H62").Value = 1950
$ws.Range("I62").Value = 1950
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1950
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1326
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 1950
$ws.Range("I65").Value = 1950
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 9750
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -6630
$ws.Range("N65").ClearContents()

$ws.Range("H74").Value = 30003
$ws.Range("I74").Value = 30003
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 30003
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -29067
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 30003
$ws.Range("I77").Value = 30003
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 150015
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -145335
$ws.Range("N77").ClearContents()

$ws.Range("H99").Value = 1110.6
$ws.Range("I99").Value = 779.5
$ws.Range("J99").Value = 1607.25
$ws.Range("K99").Value = 2338.5
$ws.Range("L99").Value = 4821.75
$ws.Range("M99").Value = -840.5
$ws.Range("N99").Value = -7817.75

$ws.Range("H129").Value = 2360.3684
$ws.Range("I129").Value = 400
$ws.Range("J129").Value = 3265.1538
$ws.Range("K129").Value = 1200
$ws.Range("L129").Value = 9795.4614
$ws.Range("M129").Value = 3800
$ws.Range("N129").Value = -19795.4614

$ws.Range("H132").Value = 1954633.2
$ws.Range("I132").Value = 2404783.2
$ws.Range("K132").Value = 7214349.600000001
$ws.Range("M132").Value = -7211819.600000001

$ws.Range("H141").Value = 1963.2245
$ws.Range("I141").Value = 1337.3903
$ws.Range("J141").Value = 5170.625
$ws.Range("K141").Value = 4012.1709
$ws.Range("L141").Value = 15511.875
$ws.Range("M141").Value = 1167.8291
$ws.Range("N141").Value = -25871.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1336.45
$ws.Range("I102").Value = 1272.1082
$ws.Range("J102").Value = 2130
$ws.Range("K102").Value = 1272.1082
$ws.Range("L102").Value = 2130
$ws.Range("M102").Value = 349.8918000000001
$ws.Range("N102").Value = -5374

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1528.6207
$ws.Range("I86").Value = 1267.0625
$ws.Range("J86").Value = 1850.5385
$ws.Range("K86").Value = 1267.0625
$ws.Range("L86").Value = 1850.5385
$ws.Range("M86").Value = -144.0625
$ws.Range("N86").Value = -4096.538500000001

$ws.Range("H89").Value = 1528.6207
$ws.Range("I89").Value = 1267.0625
$ws.Range("J89").Value = 1850.5385
$ws.Range("K89").Value = 6335.3125
$ws.Range("L89").Value = 9252.692500000001
$ws.Range("M89").Value = -719.3125
$ws.Range("N89").Value = -20484.6925

$ws.Range("H134").Value = 1918.2222
$ws.Range("I134").Value = 1125
$ws.Range("J134").Value = 2552.8
$ws.Range("K134").Value = 3375
$ws.Range("L134").Value = 7658.400000000001
$ws.Range("M134").Value = -840
$ws.Range("N134").Value = -12728.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2441.9033
$ws.Range("J31").Value = 2859.5
$ws.Range("L31").Value = 2859.5
$ws.Range("N31").Value = -3449.5

$ws.Range("H34").Value = 2441.9033
$ws.Range("J34").Value = 2859.5
$ws.Range("L34").Value = 2859.5
$ws.Range("N34").Value = -3263.5

$ws.Range("H88").Value = 24199.4
$ws.Range("I88").Value = 10311
$ws.Range("J88").Value = 27671.5
$ws.Range("K88").Value = 10311
$ws.Range("L88").Value = 27671.5
$ws.Range("M88").Value = -9905
$ws.Range("N88").Value = -28483.5

$ws.Range("H91").Value = 24199.4
$ws.Range("I91").Value = 10311
$ws.Range("J91").Value = 27671.5
$ws.Range("K91").Value = 10311
$ws.Range("L91").Value = 27671.5
$ws.Range("M91").Value = -8907
$ws.Range("N91").Value = -30479.5

$ws.Range("H99").Value = 3065.5715
$ws.Range("I99").Value = 2827.6365
$ws.Range("J99").Value = 3938
$ws.Range("K99").Value = 2827.6365
$ws.Range("L99").Value = 3938
$ws.Range("M99").Value = -1329.6365
$ws.Range("N99").Value = -6934

$ws.Range("H126").Value = 3065.5715
$ws.Range("I126").Value = 2827.6365
$ws.Range("J126").Value = 3938
$ws.Range("K126").Value = 8482.9095
$ws.Range("L126").Value = 11814
$ws.Range("M126").Value = -6012.9095
$ws.Range("N126").Value = -16754

$ws.Range("H132").Value = 1491.7142
$ws.Range("I132").Value = 1213.3226
$ws.Range("K132").Value = 3639.9678
$ws.Range("M132").Value = -1109.9678

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 391.125
$ws.Range("I113").Value = 360.84616
$ws.Range("J113").Value = 460.58823
$ws.Range("K113").Value = 1082.53848
$ws.Range("L113").Value = 1381.76469
$ws.Range("M113").Value = 1087.46152
$ws.Range("N113").Value = -5721.76469

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2638
$ws.Range("I102").Value = 2001.28
$ws.Range("J102").Value = 3475.7896
$ws.Range("K102").Value = 2001.28
$ws.Range("L102").Value = 3475.7896
$ws.Range("M102").Value = -379.28
$ws.Range("N102").Value = -6719.7896

$ws.Range("H122").Value = 1572.3148
$ws.Range("I122").Value = 1286.909
$ws.Range("J122").Value = 2020.8096
$ws.Range("K122").Value = 3860.727
$ws.Range("L122").Value = 6062.4288
$ws.Range("M122").Value = -1410.727
$ws.Range("N122").Value = -10962.4288

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2310.4443
$ws.Range("I7").Value = 2310.8
$ws.Range("J7").Value = 2310
$ws.Range("K7").Value = 2310.8
$ws.Range("L7").Value = 2310
$ws.Range("M7").Value = -2198.8
$ws.Range("N7").Value = -2534

$ws.Range("H40").Value = 2642.8572
$ws.Range("I40").Value = 2509.0908
$ws.Range("J40").Value = 3133.3333
$ws.Range("K40").Value = 2509.0908
$ws.Range("L40").Value = 3133.3333
$ws.Range("M40").Value = -2373.0908
$ws.Range("N40").Value = -3405.3333

$ws.Range("H126").Value = 2310.4443
$ws.Range("I126").Value = 2310.8
$ws.Range("J126").Value = 2310
$ws.Range("K126").Value = 6932.400000000001
$ws.Range("L126").Value = 6930
$ws.Range("M126").Value = -4462.400000000001
$ws.Range("N126").Value = -11870

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1311
$ws.Range("I126").Value = 1054.4
$ws.Range("J126").Value = 1952.5
$ws.Range("K126").Value = 3163.2
$ws.Range("L126").Value = 5857.5
$ws.Range("M126").Value = -693.2000000000003
$ws.Range("N126").Value = -10797.5

$ws.Range("H132").Value = 1229.82
$ws.Range("I132").Value = 916.2
$ws.Range("J132").Value = 1700.25
$ws.Range("K132").Value = 2748.6
$ws.Range("L132").Value = 5100.75
$ws.Range("M132").Value = -218.6000000000004
$ws.Range("N132").Value = -10160.75
